# ScriptMasterSheet.xlsx — add a "negative" variant of the testT4116 row.
#
#   Row 2 (existing testT4116 / YES row) gets highlighted green.
#   Row 3 is turned into a new "testT4116_Negative" / YES row, highlighted red
#   (previously it held a second, unrelated testT4116 / NO row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestScriptMasterSheet")

# OLE/VBA colors are packed as 0x00BBGGRR (blue high byte, red low byte),
# i.e. the reverse byte order of the usual #RRGGBB notation.
$colorGreen = 0x50D092   # fgColor FF92D050
$colorRed   = 0x0000FF   # fgColor FFFF0000

# --- Row 2: keep its values, just mark it green ---
$ws.Range("A2:C2").Interior.Color = $colorGreen

# --- Row 3: replace with the negative-test row and mark it red ---
$ws.Range("A3").Value = "testT4116_Negative"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "YES"
$ws.Range("A3:C3").Interior.Color = $colorRed

# --- Leave the selection on C3, matching the saved view state ---
$ws.Range("C3").Select()
